# Fruta / hortaliza, semanal
# Updates the weekly price records (rows 2-9) for the Caqui sheet by
# re-assigning the per-row data (Fecha, Variedad, Calidad, Volumen,
# Precio minimo/maximo/promedio, Origen, Precio $/Kg) to reflect the
# latest weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for each row (2-9), columns D, K, L, M, N, O, P, R, S
$data = @{
    2 = @{ D = 44305; K = "Mankaki"; L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 }
    3 = @{ D = 44342; K = "Mankaki"; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 }
    4 = @{ D = 45071; K = "Fuyu";    L = "Segunda"; M = 110; N = 23000; O = 24000; P = 23455; R = "Región Metropolitana"; S = 1303 }
    5 = @{ D = 44355; K = "Mankaki"; L = "Segunda"; M = 270; N = 20000; O = 21000; P = 20500; R = "Región Metropolitana"; S = 1139 }
    6 = @{ D = 44699; K = "Mankaki"; L = "Primera"; M = 250; N = 29000; O = 30000; P = 29500; R = "Región de O'Higgins"; S = 1639 }
    7 = @{ D = 45043; K = "Fuyu";    L = "Primera"; M = 300; N = 25000; O = 26000; P = 25500; R = "Región de O'Higgins"; S = 1417 }
    8 = @{ D = 44301; K = "Hachiya"; L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; R = "Región de O'Higgins"; S = 1139 }
    9 = @{ D = 44313; K = "Mankaki"; L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; R = "Región de O'Higgins"; S = 1194 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $rowData.D   # D: Fecha
    $ws.Cells.Item($row, 11).Value = $rowData.K   # K: Variedad
    $ws.Cells.Item($row, 12).Value = $rowData.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $rowData.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $rowData.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $rowData.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $rowData.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 18).Value = $rowData.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $rowData.S   # S: Precio $/Kg
}
